$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (shifts old rows 4..127 down to 5..128)
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly record.
# (Same market/region/category metadata as the record that used to sit
# in row 4 — only the date and the price columns differ.)
$ws.Range("A4").Value = 11
$ws.Range("B4").Value = "Vega Monumental Concepción"
$ws.Range("C4").Value = "Bíobío"
$ws.Range("D4").Value = 44691
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 100112043
$ws.Range("G4").Value = "Pepino ensalada"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 16000
$ws.Range("L4").Value = 17000
$ws.Range("M4").Value = 16500
$ws.Range("N4").Value = "$/caja 60 unidades"
$ws.Range("O4").Value = "Región de Arica y Parinacota"
$ws.Range("P4").Value = 275
$ws.Range("Q4").Value = 60
$ws.Range("R4").Value = "Hortaliza"
